$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp value on row 3 (slight re-run precision change)
$ws.Range("A3").Value = 45867.08370511574

# Add new row 4 with the latest automated reading
$ws.Range("A4").Value = 45867.1253028607
$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 12.72
$ws.Range("E4").Value = 89.38
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "03:00:26"

# Match the number format / style used by the other date cells in column A
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
